$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF columns, matching style of existing header cells (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-27: I column is always 1, J column mirrors H column's value
for ($row = 2; $row -le 27; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValue
}
